$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (col E) / Correspond Handback DateTime (col H)
# for rows 3 and 4 (the 83bf0102... and c03a4df5... entries) move forward in time.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-12 04:22:00"
$wsZhCn.Range("H3").Value = "2016-03-12 04:22:22"
$wsZhCn.Range("E4").Value = "2016-03-12 04:22:00"
$wsZhCn.Range("H4").Value = "2016-03-12 04:22:22"

# de-de sheet: same two rows / columns.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-12 04:22:03"
$wsDeDe.Range("H3").Value = "2016-03-12 04:22:28"
$wsDeDe.Range("E4").Value = "2016-03-12 04:22:03"
$wsDeDe.Range("H4").Value = "2016-03-12 04:22:28"
